$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1867.7667
$ws.Range("I40").Value = 1887.4286
$ws.Range("J40").Value = 1821.8889
$ws.Range("K40").Value = 1887.4286
$ws.Range("L40").Value = 1821.8889
$ws.Range("M40").Value = -1712.4286
$ws.Range("N40").Value = -2171.8889

$ws.Range("H98").Value = 2143.7856
$ws.Range("I98").Value = 431.2857
$ws.Range("K98").Value = 431.2857
$ws.Range("M98").Value = 1066.7143

$ws.Range("H122").Value = 2143.7856
$ws.Range("I122").Value = 431.2857
$ws.Range("K122").Value = 1293.8571
$ws.Range("M122").Value = 1156.1429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5629.7646
$ws.Range("I32").Value = 5356.6875
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 5356.6875
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -5069.6875
$ws.Range("N32").Value = -10573

$ws.Range("H63").Value = 27026
$ws.Range("I63").Value = 27026
$ws.Range("K63").Value = 27026
$ws.Range("M63").Value = -26340

$ws.Range("H66").Value = 27026
$ws.Range("I66").Value = 27026
$ws.Range("K66").Value = 135130
$ws.Range("M66").Value = -131698

$ws.Range("H74").Value = 1482.1666
$ws.Range("I74").Value = 1278.6
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 1278.6
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -404.5999999999999
$ws.Range("N74").Value = -4248

$ws.Range("H77").Value = 1482.1666
$ws.Range("I77").Value = 1278.6
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 6393
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -2025
$ws.Range("N77").Value = -21236

$ws.Range("H122").Value = 1631.25
$ws.Range("I122").Value = 930
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 2790
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -340
$ws.Range("N122").Value = -13300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 299999
$ws.Range("J42").Value = 299999
$ws.Range("L42").Value = 299999
$ws.Range("N42").Value = -300655

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0

$ws.Range("H99").Value = 2398.6667
$ws.Range("I99").Value = 2848.5
$ws.Range("J99").Value = 1499
$ws.Range("K99").Value = 2848.5
$ws.Range("L99").Value = 1499
$ws.Range("M99").Value = -1350.5
$ws.Range("N99").Value = -4495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 978484
$ws.Range("I6").Value = 1412754.6
$ws.Range("J6").Value = 1375
$ws.Range("K6").Value = 1412754.6
$ws.Range("L6").Value = 1375
$ws.Range("M6").Value = -1412641.6
$ws.Range("N6").Value = -1601

$ws.Range("H19").Value = 431.6
$ws.Range("I19").Value = 39.5
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 39.5
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 130.5
$ws.Range("N19").Value = -2340

$ws.Range("H22").Value = 449.75
$ws.Range("I22").Value = 466.33334
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 466.33334
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -116.33334
$ws.Range("N22").Value = -1100

$ws.Range("H24").Value = 431.6
$ws.Range("I24").Value = 39.5
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 39.5
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 130.5
$ws.Range("N24").Value = -2340

$ws.Range("H31").Value = 2143.2856
$ws.Range("I31").Value = 1640.5
$ws.Range("K31").Value = 1640.5
$ws.Range("M31").Value = -1345.5

$ws.Range("H32").Value = 1799.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1799.5
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = 1799.5
$ws.Range("N32").Value = -2431.5

$ws.Range("H34").Value = 2143.2856
$ws.Range("I34").Value = 1640.5
$ws.Range("K34").Value = 1640.5
$ws.Range("M34").Value = -1438.5

$ws.Range("H92").Value = 14600
$ws.Range("J92").Value = 14600
$ws.Range("L92").Value = 14600
$ws.Range("N92").Value = -19592

$ws.Range("H122").Value = 3290.875
$ws.Range("I122").Value = 3290.875
$ws.Range("K122").Value = 9872.625
$ws.Range("M122").Value = -7422.625

$ws.Range("H134").Value = 1174.875
$ws.Range("I134").Value = 1018.1667
$ws.Range("J134").Value = 1645
$ws.Range("K134").Value = 3054.5001
$ws.Range("L134").Value = 4935
$ws.Range("M134").Value = -519.5001000000002
$ws.Range("N134").Value = -10005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7999
$ws.Range("I3").Value = 7999
$ws.Range("K3").Value = 23997
$ws.Range("M3").Value = -23885

$ws.Range("H114").Value = 4139
$ws.Range("I114").Value = 4030.5
$ws.Range("J114").Value = 4225.8
$ws.Range("K114").Value = 12091.5
$ws.Range("L114").Value = 12677.4
$ws.Range("M114").Value = -8837.5
$ws.Range("N114").Value = -19185.4

$ws.Range("H119").Value = 3450
$ws.Range("I119").Value = 3450
$ws.Range("K119").Value = 10350
$ws.Range("M119").Value = -5512

$ws.Range("H130").Value = 1892
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 1892
$ws.Range("K130").Value = 0
$ws.Range("M130").Value = 5676
$ws.Range("N130").Value = -15716

$ws.Range("H131").Value = 296959.16
$ws.Range("J131").Value = 347777.9
$ws.Range("L131").Value = 1043333.7
$ws.Range("N131").Value = -1053413.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 26144
$ws.Range("J39").Value = 26144
$ws.Range("L39").Value = 26144
$ws.Range("N39").Value = -27208

$ws.Range("H102").Value = 1899.125
$ws.Range("I102").Value = 1899.125
$ws.Range("K102").Value = 1899.125
$ws.Range("M102").Value = -277.125

$ws.Range("H122").Value = 2092
$ws.Range("I122").Value = 2244.2
$ws.Range("J122").Value = 1838.3334
$ws.Range("K122").Value = 6732.599999999999
$ws.Range("L122").Value = 5515.0002
$ws.Range("M122").Value = -4282.599999999999
$ws.Range("N122").Value = -10415.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 100
$ws.Range("K17").Value = 100
$ws.Range("M17").Value = 70

$ws.Range("H22").Value = 3309.8
$ws.Range("I22").Value = 2850
$ws.Range("J22").Value = 3999.5
$ws.Range("K22").Value = 2850
$ws.Range("L22").Value = 3999.5
$ws.Range("M22").Value = -2555
$ws.Range("N22").Value = -4589.5

$ws.Range("H27").Value = 3309.8
$ws.Range("I27").Value = 2850
$ws.Range("J27").Value = 3999.5
$ws.Range("K27").Value = 2850
$ws.Range("L27").Value = 3999.5
$ws.Range("M27").Value = -2743
$ws.Range("N27").Value = -4213.5

$ws.Range("H32").Value = 4999.5
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 4999
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 4999
$ws.Range("M32").Value = -4683
$ws.Range("N32").Value = -5633

$ws.Range("H55").Value = 178.16667
$ws.Range("I55").Value = 92.625
$ws.Range("J55").Value = 349.25
$ws.Range("K55").Value = 92.625
$ws.Range("L55").Value = 349.25
$ws.Range("M55").Value = 80.375
$ws.Range("N55").Value = -695.25

$ws.Range("H61").Value = 1902.9
$ws.Range("I61").Value = 931.6667
$ws.Range("J61").Value = 2319.1428
$ws.Range("K61").Value = 931.6667
$ws.Range("L61").Value = 2319.1428
$ws.Range("M61").Value = -729.6667
$ws.Range("N61").Value = -2723.1428

$ws.Range("H113").Value = 1902.9
$ws.Range("I113").Value = 931.6667
$ws.Range("J113").Value = 2319.1428
$ws.Range("K113").Value = 931.6667
$ws.Range("L113").Value = 2319.1428
$ws.Range("M113").Value = 1238.3333
$ws.Range("N113").Value = -6659.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 44998
$ws.Range("J54").Value = 44998
$ws.Range("L54").Value = 44998
$ws.Range("N54").Value = -46038

$ws.Range("H55").Value = 9408.166999999999
$ws.Range("I55").Value = 5023.25
$ws.Range("K55").Value = 5023.25
$ws.Range("M55").Value = -4746.25

$ws.Range("H68").Value = 81028.336
$ws.Range("J68").Value = 81028.336
$ws.Range("L68").Value = 81028.336
$ws.Range("N68").Value = -82650.336

$ws.Range("H71").Value = 81028.336
$ws.Range("J71").Value = 81028.336
$ws.Range("L71").Value = 243085.008
$ws.Range("N71").Value = -251197.008

$ws.Range("H100").Value = 5556758.5
$ws.Range("I100").Value = 6251166
$ws.Range("K100").Value = 12502332
$ws.Range("M100").Value = -12501791
